# Add a second worksheet "Data2" (placed after the existing "Data" sheet)
# that holds a header row used by the new "verify text and color" test
# helper: Name / Address / Contact / City / PIN Code.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("Data")

# Insert the new sheet right after "Data" so tab order becomes Data, Data2.
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "Data2"

$ws.Range("A1").Value = "Name "
$ws.Range("B1").Value = "Address "
$ws.Range("C1").Value = "Contact"
$ws.Range("D1").Value = "City"
$ws.Range("E1").Value = "PIN Code"

# Leave the selection on the last header cell, matching the saved state.
$ws.Range("E1").Select() | Out-Null
